# Report.xlsx update: correct D13 ("Time") and log the 03.01.2018 entry as row 14
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D13 previously held "8.30 to 4.30" ; correct it to "8.30 to 4.45"
$ws.Range("D13").Value = "8.30 to 4.45"

# Copy row 13's cell formatting down onto the new row 14 before filling it in,
# so the new entries match the table's existing look (font/fill/alignment)
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Append the new day's entry as row 14
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Poongodi.R"

# "03.01.2018" parses as a valid date, so a plain Value assignment would get
# auto-converted to a date serial on entry. Enter it as a quoted-text formula
# (never subject to date autodetection) and then flatten it down to a plain
# value in place, so the cell ends up holding literal text, like the rest of
# the Date column, without disturbing its existing number format/style.
$ws.Range("C14").Formula = '="03.01.2018"'
$ws.Range("C14").Copy()
$ws.Range("C14").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("D14").Value = "8.30 to 4.45"
$ws.Range("E14").Value = "PHP select,insert,update queries, attend chatbot class"
$ws.Range("F14").Value = "completed"

# Move the active selection, matching the cursor position after the edit
$ws.Range("E20").Select()
